$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.744.04'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '2.097.55'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.52'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.18'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.89%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.86'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +7.57%  '
$ws.Range('D13').Value = '2.408.11'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.14'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.807'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.11%  '
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '2.094.47'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '38.748.04'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.93'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.06'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.46'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.35'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.70%  '
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.86'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.57'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('E28').Value = '  +6.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.42'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.52%  '
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.77'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('E35').Value = '  +2.55%  '
$ws.Range('E36').Value = '  +3.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.41'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.56'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.22'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.44%  '
$ws.Range('D43').Value = '1.534.67'
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.81'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.98%  '
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.12'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.97'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').Value = '2.291.71'
$ws.Range('E51').Value = '  -0.09%  '
